$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 27
$ws.Range("D2").Value = 44208
$ws.Range("K2").Value = "Rich Lady"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 29000
$ws.Range("P2").Value = 28500
$ws.Range("Q2").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1583
$ws.Range("T2").Value = 18

# Row 3 <- old row 23
$ws.Range("D3").Value = 44160
$ws.Range("K3").Value = "Early Majestic"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("Q3").Value = "`$/caja 18 kilos granel"
$ws.Range("R3").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S3").Value = 1361
$ws.Range("T3").Value = 18

# Row 4 <- old row 22
$ws.Range("D4").Value = 44201
$ws.Range("K4").Value = "Flavor Crest"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 23000
$ws.Range("P4").Value = 22500
$ws.Range("Q4").Value = "`$/caja 18 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1250
$ws.Range("T4").Value = 18

# Row 5 <- old row 12
$ws.Range("D5").Value = 44216
$ws.Range("K5").Value = "Andross"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17500
$ws.Range("Q5").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R5").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S5").Value = 1094
$ws.Range("T5").Value = 16

# Row 6 <- old row 7
$ws.Range("D6").Value = 44229
$ws.Range("K6").Value = "Doctor Davis"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 320
$ws.Range("N6").Value = 19000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 19500
$ws.Range("Q6").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1083
$ws.Range("T6").Value = 18

# Row 7 <- old row 2
$ws.Range("D7").Value = 44203
$ws.Range("K7").Value = "Carson"
$ws.Range("L7").Value = "Tercera"
$ws.Range("M7").Value = 270
$ws.Range("N7").Value = 19000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19500
$ws.Range("Q7").Value = "`$/caja 18 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1083
$ws.Range("T7").Value = 18

# Row 8 <- old row 3
$ws.Range("D8").Value = 44203
$ws.Range("K8").Value = "Flavor Crest"
$ws.Range("L8").Value = "Tercera"
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 17000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 17500
$ws.Range("Q8").Value = "`$/caja 18 kilos granel"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 972
$ws.Range("T8").Value = 18

# Row 9 <- old row 13
$ws.Range("D9").Value = 44258
$ws.Range("K9").Value = "Doctor Davis"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 9500
$ws.Range("Q9").Value = "`$/bandeja 10 kilos empedrada"
$ws.Range("R9").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S9").Value = 950
$ws.Range("T9").Value = 10

# Row 10 <- old row 11
$ws.Range("D10").Value = 44174
$ws.Range("K10").Value = "Kurakata"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 21000
$ws.Range("P10").Value = 20500
$ws.Range("Q10").Value = "`$/caja 18 kilos granel"
$ws.Range("R10").Value = "Región de Coquimbo"
$ws.Range("S10").Value = 1139
$ws.Range("T10").Value = 18

# Row 11 <- old row 9
$ws.Range("D11").Value = 44236
$ws.Range("K11").Value = "Doctor Davis"
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20500
$ws.Range("Q11").Value = "`$/caja 18 kilos granel"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 1139
$ws.Range("T11").Value = 18

# Row 12 <- old row 18
$ws.Range("D12").Value = 44209
$ws.Range("K12").Value = "Carson"
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 19000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19500
$ws.Range("Q12").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 1083
$ws.Range("T12").Value = 18

# Row 13 <- old row 19
$ws.Range("D13").Value = 44209
$ws.Range("K13").Value = "Royal Glory"
$ws.Range("L13").Value = "Tercera"
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 19000
$ws.Range("P13").Value = 18500
$ws.Range("Q13").Value = "`$/caja 18 kilos granel"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 1028
$ws.Range("T13").Value = 18

# Row 14 <- old row 20
$ws.Range("D14").Value = 44278
$ws.Range("K14").Value = "Phillips Cling"
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 250
$ws.Range("N14").Value = 23000
$ws.Range("O14").Value = 24000
$ws.Range("P14").Value = 23500
$ws.Range("Q14").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R14").Value = "Región de O'Higgins"
$ws.Range("S14").Value = 1306
$ws.Range("T14").Value = 18

# Row 15 <- old row 24
$ws.Range("D15").Value = 44217
$ws.Range("K15").Value = "Flavor Crest"
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 19000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 19500
$ws.Range("Q15").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 1083
$ws.Range("T15").Value = 18

# Row 16 <- old row 25
$ws.Range("D16").Value = 44217
$ws.Range("K16").Value = "Royal Glory"
$ws.Range("L16").Value = "Tercera"
$ws.Range("M16").Value = 270
$ws.Range("N16").Value = 17000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 17500
$ws.Range("Q16").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R16").Value = "Región de O'Higgins"
$ws.Range("S16").Value = 972
$ws.Range("T16").Value = 18

# Row 17 <- old row 10
$ws.Range("D17").Value = 44215
$ws.Range("K17").Value = "Andross"
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 19000
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 19500
$ws.Range("Q17").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R17").Value = "Región de O'Higgins"
$ws.Range("S17").Value = 1083
$ws.Range("T17").Value = 18

# Row 18 <- old row 14
$ws.Range("D18").Value = 44222
$ws.Range("K18").Value = "Elegant Lady"
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 21000
$ws.Range("P18").Value = 20500
$ws.Range("Q18").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R18").Value = "Región de O'Higgins"
$ws.Range("S18").Value = 1139
$ws.Range("T18").Value = 18

# Row 19 <- old row 4
$ws.Range("D19").Value = 44210
$ws.Range("K19").Value = "Carson"
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 19000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 19500
$ws.Range("Q19").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 1083
$ws.Range("T19").Value = 18

# Row 20 <- old row 5
$ws.Range("D20").Value = 44210
$ws.Range("K20").Value = "Rich Lady"
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 270
$ws.Range("N20").Value = 19000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 19500
$ws.Range("Q20").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 1083
$ws.Range("T20").Value = 18

# Row 21 <- old row 6
$ws.Range("D21").Value = 44210
$ws.Range("K21").Value = "Royal Glory"
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 300
$ws.Range("N21").Value = 19000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 19500
$ws.Range("Q21").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 1083
$ws.Range("T21").Value = 18

# Row 22 <- old row 17
$ws.Range("D22").Value = 44257
$ws.Range("K22").Value = "September Sweet"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 19000
$ws.Range("O22").Value = 20000
$ws.Range("P22").Value = 19500
$ws.Range("Q22").Value = "`$/caja 18 kilos granel"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 1083
$ws.Range("T22").Value = 18

# Row 23 <- old row 15
$ws.Range("D23").Value = 44273
$ws.Range("K23").Value = "Doctor Davis"
$ws.Range("L23").Value = "Segunda"
$ws.Range("M23").Value = 280
$ws.Range("N23").Value = 22000
$ws.Range("O23").Value = 23000
$ws.Range("P23").Value = 22500
$ws.Range("Q23").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R23").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S23").Value = 1250
$ws.Range("T23").Value = 18

# Row 24 <- old row 16
$ws.Range("D24").Value = 44273
$ws.Range("K24").Value = "Phillips Cling"
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 250
$ws.Range("N24").Value = 22000
$ws.Range("O24").Value = 23000
$ws.Range("P24").Value = 22500
$ws.Range("Q24").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 1250
$ws.Range("T24").Value = 18

# Row 25 <- old row 21
$ws.Range("D25").Value = 44161
$ws.Range("K25").Value = "Florida King"
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 250
$ws.Range("N25").Value = 24000
$ws.Range("O25").Value = 25000
$ws.Range("P25").Value = 24500
$ws.Range("Q25").Value = "`$/caja 18 kilos granel"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 1361
$ws.Range("T25").Value = 18

# Row 27 <- old row 8
$ws.Range("D27").Value = 44167
$ws.Range("K27").Value = "Florida King"
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 250
$ws.Range("N27").Value = 25000
$ws.Range("O27").Value = 26000
$ws.Range("P27").Value = 25500
$ws.Range("Q27").Value = "`$/caja 18 kilos granel"
$ws.Range("R27").Value = "Región de O'Higgins"
$ws.Range("S27").Value = 1417
$ws.Range("T27").Value = 18
